# Fix de imágenes de los paretos
# Insert a new "Metodo" column at the front of the decision table on Hoja1,
# shifting the existing Rx/Ry/CL/Entropia/SSIM columns one place to the right,
# and populate it with the method name for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Shift existing data (A:E) one column to the right, to make room for the
# new "Metodo" column at A.
$ws.Columns.Item(1).Insert()

# Headers: "Metodo" is brand new in column A; B1:F1 keep the original
# Rx/Ry/CL/Entropia/SSIM labels (now shifted one column right).
$ws.Range("A1").Value = "Metodo"
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

$metodos = @("SMARTER", "Fuzzy", "TOPSIS", "GRA", "CODAS", "MABAC", "VIKOR", "PROMETHEE II")
for ($i = 0; $i -lt $metodos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $metodos[$i]
}

# Re-apply best-fit column widths for the new/changed columns (Rx/Ry/CL/
# Entropia/SSIM in D:F already carry their original best-fit width and are
# left untouched).
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
